$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 was missing the PriceChange / UpDown classification - fill it in now
# (the "moved against" check was producing false positives, so this was re-run)
$ws.Range("X3").Value = 1.0200049999999976
$ws.Range("Y3").Value = "Up"

# Append the next day's prediction as a new row (row 4)
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 42641.892604166664

$ws.Range("B4").Value = -15
$ws.Range("C4").Value = "Strong Sell"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = "Random"
$ws.Range("Q4").Value = 52.89259217263573
$ws.Range("R4").Value = 0.85

$ws.Range("S3").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 0.0202

$ws.Range("T3").Copy()
$ws.Range("T4").PasteSpecial(-4122)
$ws.Range("T4").Value = -0.0172

$ws.Range("U4").Value = 15
$ws.Range("V4").Value = "N/A"
$ws.Range("W4").Value = 0

# "Strong Sell" is now the longest value in the Verdict column, so it no
# longer fits the old best-fit width - widen column C to suit.
$ws.Columns.Item(3).ColumnWidth = 8
